$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.850.11'
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').Value = '1.876.19'
$ws.Range('E3').Value = '  +0.12%  '
$ws.Range('D4').Value = '''0.9996'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''0.7290'
$ws.Range('E5').Value = '  -1.28%  '
$ws.Range('D6').Value = '''241.96'
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('D7').Value = '''1.000'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '''0.3145'
$ws.Range('E8').Value = '  -0.28%  '
$ws.Range('D9').Value = '''0.07552'
$ws.Range('E9').Value = '  +5.60%  '
$ws.Range('D10').Value = '''24.65'
$ws.Range('E10').Value = '  -0.12%  '
$ws.Range('D11').Value = '''0.08180'
$ws.Range('E11').Value = '  -2.73%  '
$ws.Range('D12').Value = '''0.7474'
$ws.Range('E12').Value = '  -0.40%  '
$ws.Range('D13').Value = '''5.365'
$ws.Range('E13').Value = '  -1.01%  '
$ws.Range('D14').Value = '1.875.02'
$ws.Range('E14').Value = '  -0.63%  '
$ws.Range('D15').Value = '''92.79'
$ws.Range('D16').Value = '29.817.78'
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('D17').Value = '''6.020'
$ws.Range('E17').Value = '  -1.25%  '
$ws.Range('D18').Value = '''248.15'
$ws.Range('E18').Value = '  +2.17%  '
$ws.Range('D19').Value = '''13.48'
$ws.Range('E19').Value = '  -0.75%  '
$ws.Range('D20').Value = '''0.000007907'
$ws.Range('E20').Value = '  +1.21%  '
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').Value = '2.118.69'
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').Value = '''0.9999'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').Value = '''7.754'
$ws.Range('E24').Value = '  -2.94%  '
$ws.Range('E25').Value = '  -1.55%  '
$ws.Range('D26').Value = '''9.271'
$ws.Range('E26').Value = '  -0.16%  '
$ws.Range('D27').Value = '''164.40'
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('D28').Value = '''18.63'
$ws.Range('E28').Value = '  +0.13%  '
$ws.Range('D29').Value = '''2.016'
$ws.Range('E29').Value = '  -0.94%  '
$ws.Range('D30').Value = '''1.447'
$ws.Range('E30').Value = '  -2.97%  '
$ws.Range('D31').Value = '''4.543'
$ws.Range('E31').Value = '  -1.08%  '
$ws.Range('D32').Value = '''1.526'
$ws.Range('E32').Value = '  -0.15%  '
$ws.Range('D33').Value = '''4.191'
$ws.Range('E33').Value = '  -1.29%  '
$ws.Range('D34').Value = '''0.05417'
$ws.Range('E34').Value = '  +1.84%  '
$ws.Range('D35').Value = '''1.234'
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('D36').Value = '''0.7435'
$ws.Range('E36').Value = '  -1.40%  '
$ws.Range('D37').Value = '''1.004'
$ws.Range('E37').Value = '  +0.85%  '
$ws.Range('D38').Value = '''2.705'
$ws.Range('E38').Value = '  +0.36%  '
$ws.Range('D39').Value = '''0.01932'
$ws.Range('E39').Value = '  -0.82%  '
$ws.Range('D40').Value = '''2.744'
$ws.Range('E40').Value = '  -0.34%  '
$ws.Range('D41').Value = '''0.4479'
$ws.Range('E41').Value = '  -0.67%  '
$ws.Range('D42').Value = '''0.8889'
$ws.Range('E42').Value = '  +3.84%  '
$ws.Range('D43').Value = '''6.002'
$ws.Range('E43').Value = '  -0.81%  '
$ws.Range('D44').Value = '''71.81'
$ws.Range('E44').Value = '  -0.51%  '
$ws.Range('D45').Value = '''104.43'
$ws.Range('E45').Value = '  +1.32%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '1.042.79'
$ws.Range('E46').Value = '  -6.01%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').Value = '''1.001'
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').Value = '''7.498'
$ws.Range('E48').Value = '  -1.95%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '''1.826'
$ws.Range('E49').Value = '  -0.72%  '
$ws.Range('D50').Value = '''9.690'
$ws.Range('E50').Value = '  +1.69%  '
$ws.Range('D51').Value = '2.019.59'
$ws.Range('E51').Value = '  +0.27%  '
